$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "رقم البوليصه" (policy/police id) column entirely, shifting
# everything else one column to the left.
$ws.Columns.Item(1).Delete()

# Update selection to mirror the post-edit workbook (whole column A selected).
$ws.Range("A1:A1048576").Select()
